# feat: add 2022-Q3 data
#
# 1. Insert a brand-new "2022-Q3" worksheet right before the existing
#    "2022-Q2" sheet (so tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1)
#    and populate it with the Q3 fund-holding table.
# 2. Insert a new row 2 at the top of the "总计" (totals) sheet's data and
#    fill it in with the Q3 summary figures, pushing the existing Q2/Q1
#    rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: new "2022-Q3" sheet
# ---------------------------------------------------------------------------
$ws2022Q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($ws2022Q2)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3hdr = $q3.Range("B1:H1")
$q3hdr.Font.Bold = $true
$q3hdr.HorizontalAlignment = -4108
$q3hdr.VerticalAlignment = -4160
$q3hdr.Borders.LineStyle = 1

# Fund code / name / percentage-ish columns are stored as text in the
# source workbook (e.g. "007355" keeps its leading zero) - force the "@"
# text format before writing so COM doesn't silently coerce them to numbers.
$q3.Range("B2:G6").NumberFormat = "@"

$q3data = @(
  @(0, "506006", "汇添富科创板2年定期开放混合", "16.61", "88.62", "4.25", "0.7059", 8),
  @(1, "007355", "汇添富科技创新灵活配置混合A", "14.72", "91.95", "3.06", "0.4504", 10),
  @(2, "007356", "汇添富科技创新灵活配置混合C", "2.53", "91.95", "3.06", "0.0774", 10),
  @(3, "014232", "博时专精特新主题混合A", "3.14", "81.93", "1.66", "0.0521", 5),
  @(4, "014233", "博时专精特新主题混合C", "2.69", "81.93", "1.66", "0.0447", 5)
)

$r = 2
foreach ($row in $q3data) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Drop the temporary text format now the values are locked in as text, so
# the data cells end up with no explicit style - matching the sibling
# "2022-Q2" / "2022-Q1" sheets.
$q3.Range("B2:G6").ClearFormats()

# Row-number column keeps the bold/bordered/centered look used elsewhere.
$q3acol = $q3.Range("A2:A6")
$q3acol.Font.Bold = $true
$q3acol.HorizontalAlignment = -4108
$q3acol.VerticalAlignment = -4160
$q3acol.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Step 2: update the "总计" (totals) sheet with the new Q3 row
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Push the existing rows (2022-Q2 / 2022-Q1) down and open up row 2.
$totals.Rows.Item(2).Insert()

# Row 2 lost its own formatting on insert but also doesn't carry the
# row-label style automatically - copy it from row 3 (still intact).
$totals.Range("A2:D2").ClearFormats()
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)
$totals.Application.CutCopyMode = $false

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 1.33

# Renumber the row-index column for the rows that shifted down.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2

Write-Output "done"
